$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: PT Prima Tunggal Mandiri (Shell) - updated lease dates, income, and payment scheme
$ws.Range("B8").Value = 45939
$ws.Range("C8").Value = 46304
$ws.Range("D8").Value = ""
$ws.Range("F8").Value = 170000000
$ws.Range("G8").Value = 170000000
$ws.Range("H8").Value = "Full Lease Upfront"

# Row 15: PT Mandiri Utama Finance - projected income now matches actual income
$ws.Range("F15").Value = 280000000
